# Updates the cryptos price list (columns D = Price, E = Volume(1h)) for rows 2-51
# to the latest values, as captured by the GitHub Actions scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is: cell address -> new text value.
# New values are written with a leading apostrophe so Excel keeps them as text
# (matching the original "Price"/"Volume(1h)" columns, which are plain text),
# and the cell style is reset to "Normal" afterwards so no stray number format sticks.
$updates = [ordered]@{
    "D2" = '''29.125.36'
    "E2" = '''  +1.32%  '
    "D3" = '''1.903.96'
    "E3" = '''  +1.61%  '
    "D4" = '''1.003'
    "E4" = '''  -0.18%  '
    "D5" = '''327.00'
    "E5" = '''  +0.82%  '
    "E6" = '''  -0.10%  '
    "D7" = '''0.4606'
    "E7" = '''  +0.02%  '
    "D8" = '''0.3930'
    "E8" = '''  +1.50%  '
    "D9" = '''46.83'
    "E9" = '''  +1.42%  '
    "D10" = '''0.07931'
    "E10" = '''  +0.92%  '
    "D11" = '''0.9993'
    "E11" = '''  +1.11%  '
    "E12" = '''  +2.04%  '
    "D13" = '''1.921.83'
    "E13" = '''  +2.70%  '
    "D14" = '''7.066'
    "E14" = '''  +1.09%  '
    "D15" = '''5.758'
    "D16" = '''0.06949'
    "E16" = '''  -0.44%  '
    "D17" = '''88.32'
    "E17" = '''  -0.14%  '
    "D18" = '''1.003'
    "E18" = '''  -0.13%  '
    "E19" = '''  +0.20%  '
    "D20" = '''17.11'
    "E20" = '''  +1.93%  '
    "E21" = '''  -0.14%  '
    "D22" = '''29.136.95'
    "E22" = '''  +1.36%  '
    "D23" = '''5.356'
    "E23" = '''  +1.41%  '
    "E24" = '''  +0.31%  '
    "D25" = '''2.105.00'
    "E25" = '''  +0.39%  '
    "D26" = '''2.054'
    "E26" = '''  -2.16%  '
    "D27" = '''156.49'
    "E27" = '''  +2.33%  '
    "D28" = '''19.43'
    "E28" = '''  +0.83%  '
    "D29" = '''6.125'
    "E29" = '''  +4.45%  '
    "D30" = '''1.993'
    "E30" = '''  +0.86%  '
    "D31" = '''118.72'
    "E31" = '''  -0.12%  '
    "D32" = '''0.09373'
    "E32" = '''  +0.58%  '
    "D33" = '''0.9275'
    "E33" = '''  +0.79%  '
    "E34" = '''  +0.28%  '
    "D35" = '''1.348'
    "E35" = '''  +0.57%  '
    "D36" = '''3.269'
    "E36" = '''  -1.57%  '
    "D37" = '''1.201'
    "E37" = '''  +4.44%  '
    "D38" = '''0.05826'
    "E38" = '''  +0.88%  '
    "D39" = '''0.02103'
    "E39" = '''  +1.55%  '
    "D40" = '''7.921'
    "E40" = '''  +3.21%  '
    "D41" = '''1.001'
    "E41" = '''  -0.21%  '
    "D42" = '''0.5738'
    "E42" = '''  +1.68%  '
    "D43" = '''0.1798'
    "E43" = '''  +0.65%  '
    "D44" = '''9.937'
    "E44" = '''  +1.02%  '
    "D45" = '''11.98'
    "E45" = '''  +1.46%  '
    "E46" = '''  +4.70%  '
    "D47" = '''0.5409'
    "E47" = '''  +2.13%  '
    "D48" = '''0.07073'
    "E48" = '''  -1.91%  '
    "E49" = '''  +2.48%  '
    "D50" = '''2.553'
    "E50" = '''  +5.62%  '
    "D51" = '''113.13'
    "E51" = '''  -0.31%  '
}

foreach ($addr in $updates.Keys) {
    $c = $ws.Range($addr)
    $c.Value = $updates[$addr]
    $c.Style = "Normal"
}
